$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# New header cells (bold/bordered/centered header style matches the rest of row 1)
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# Team record values for every data row (2-44)
for ($row = 2; $row -le 44; $row++) {
    $ws.Cells.Item($row, 30).Value = 71
    $ws.Cells.Item($row, 31).Value = 91
    $ws.Cells.Item($row, 32).Value = 0
}
